# Remove the "Agenda" slide from the presentation.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isAgenda = $false

    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text.Trim() -eq "Agenda") {
                $isAgenda = $true
                break
            }
        }
    }

    if ($isAgenda) {
        $slide.Delete()
    }
}
